# Auto-generated Excel COM-interop script to apply the described edit
$wb = $excel.ActiveWorkbook

# --- 1. Metadata sheet: update "Last Updated" timestamp ---
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("A2").Value = "05 Nov 2025, 03:15 PM"

# --- 2. Industry Analysis sheet: update "1 Year" (column F) values for rows 2-76 ---
$wsIndustry = $wb.Worksheets.Item("Industry Analysis")
$wsIndustry.Range("F2").Value = 21.0016
$wsIndustry.Range("F3").Value = -16.2396
$wsIndustry.Range("F4").Value = 27.1317
$wsIndustry.Range("F5").Value = -50.6494
$wsIndustry.Range("F6").Value = 53.2813
$wsIndustry.Range("F7").Value = -8.1062
$wsIndustry.Range("F8").Value = -9.5521
$wsIndustry.Range("F9").Value = 36.3756
$wsIndustry.Range("F10").Value = -6.1314
$wsIndustry.Range("F11").Value = 31.9081
$wsIndustry.Range("F12").Value = -18.4955
$wsIndustry.Range("F13").Value = 14.0155
$wsIndustry.Range("F14").Value = -36.0718
$wsIndustry.Range("F15").Value = -0.1622
$wsIndustry.Range("F16").Value = 0.1459
$wsIndustry.Range("F17").Value = -22.0012
$wsIndustry.Range("F18").Value = 1.0561
$wsIndustry.Range("F19").Value = -27.708
$wsIndustry.Range("F20").Value = 47.7309
$wsIndustry.Range("F21").Value = 12.0959
$wsIndustry.Range("F22").Value = 95.1491
$wsIndustry.Range("F23").Value = -50.2657
$wsIndustry.Range("F24").Value = -13.3427
$wsIndustry.Range("F25").Value = -9.9316
$wsIndustry.Range("F26").Value = 5.8244
$wsIndustry.Range("F27").Value = -32.7692
$wsIndustry.Range("F28").Value = -24.8224
$wsIndustry.Range("F29").Value = -18.4191
$wsIndustry.Range("F30").Value = 25.8569
$wsIndustry.Range("F31").Value = 58.4712
$wsIndustry.Range("F32").Value = -3.3862
$wsIndustry.Range("F33").Value = -6.3282
$wsIndustry.Range("F34").Value = 27.7203
$wsIndustry.Range("F35").Value = 4.4873
$wsIndustry.Range("F36").Value = -4.9458
$wsIndustry.Range("F37").Value = 3.6074
$wsIndustry.Range("F38").Value = -23.3973
$wsIndustry.Range("F39").Value = 8.7355
$wsIndustry.Range("F40").Value = -5.8541
$wsIndustry.Range("F41").Value = -8.3934
$wsIndustry.Range("F42").Value = 20.3818
$wsIndustry.Range("F43").Value = 14.3164
$wsIndustry.Range("F44").Value = -12.6846
$wsIndustry.Range("F45").Value = 28.4075
$wsIndustry.Range("F46").Value = -1.1135
$wsIndustry.Range("F47").Value = -37.1997
$wsIndustry.Range("F48").Value = -29.8569
$wsIndustry.Range("F49").Value = -27.5511
$wsIndustry.Range("F50").Value = -49.7478
$wsIndustry.Range("F51").Value = -51.8002
$wsIndustry.Range("F52").Value = -38.5254
$wsIndustry.Range("F53").Value = -12.4886
$wsIndustry.Range("F54").Value = -5.0725
$wsIndustry.Range("F55").Value = -17.7445
$wsIndustry.Range("F56").Value = -26.636
$wsIndustry.Range("F57").Value = -29.3361
$wsIndustry.Range("F58").Value = -11.9574
$wsIndustry.Range("F59").Value = -24.5687
$wsIndustry.Range("F60").Value = -12.3
$wsIndustry.Range("F61").Value = -10.9446
$wsIndustry.Range("F62").Value = -17.1229
$wsIndustry.Range("F63").Value = -9.5038
$wsIndustry.Range("F64").Value = 54.2749
$wsIndustry.Range("F65").Value = -43.4736
$wsIndustry.Range("F66").Value = 13.2687
$wsIndustry.Range("F67").Value = 12.7149
$wsIndustry.Range("F68").Value = 24.8057
$wsIndustry.Range("F69").Value = -17.0328
$wsIndustry.Range("F70").Value = -6.8927
$wsIndustry.Range("F71").Value = 13.6034
$wsIndustry.Range("F72").Value = 3.9995
$wsIndustry.Range("F73").Value = -16.226
$wsIndustry.Range("F74").Value = -16.2448
$wsIndustry.Range("F75").Value = 28.6924
$wsIndustry.Range("F76").Value = 48.9752

# --- 3. Stock List sheet: rows shift up by one (oldest entry dropped, newest appended) ---
$wsStock = $wb.Worksheets.Item("Stock List")
$wsStock.Range("B2").Value = "NIFTYCASE"
$wsStock.Range("C2").Value = "NIFTYCASE"
$wsStock.Range("D2").Value = 10.19
$wsStock.Range("E2").Value = -0.5854
$wsStock.Range("H2").Value = 0
$wsStock.Range("B3").Value = "MOMENTUM30"
$wsStock.Range("C3").Value = "MOMENTUM30"
$wsStock.Range("D3").Value = 31.54
$wsStock.Range("E3").Value = -0.6614
$wsStock.Range("H3").Value = 0
$wsStock.Range("B4").Value = "CANHLIFE"
$wsStock.Range("C4").Value = "CANHLIFE"
$wsStock.Range("D4").Value = 118.46
$wsStock.Range("E4").Value = 0.6286
$wsStock.Range("H4").Value = 11253.7
$wsStock.Range("B5").Value = "FLEXIADD"
$wsStock.Range("C5").Value = "FLEXIADD"
$wsStock.Range("D5").Value = 10.64
$wsStock.Range("E5").Value = -1.0233
$wsStock.Range("H5").Value = 0
$wsStock.Range("B6").Value = "MOENERGY"
$wsStock.Range("C6").Value = "MOENERGY"
$wsStock.Range("D6").Value = 36.3
$wsStock.Range("E6").Value = -0.6568000000000001
$wsStock.Range("H6").Value = 0
$wsStock.Range("B7").Value = "MONIFTY100"
$wsStock.Range("C7").Value = "MONIFTY100"
$wsStock.Range("D7").Value = 26.49
$wsStock.Range("E7").Value = 0.3409
$wsStock.Range("H7").Value = 0
$wsStock.Range("B8").Value = "RUBICON"
$wsStock.Range("C8").Value = "RUBICON"
$wsStock.Range("D8").Value = 652.65
$wsStock.Range("E8").Value = -0.1453
$wsStock.Range("H8").Value = 10752.4289
$wsStock.Range("B9").Value = "CRAMC"
$wsStock.Range("C9").Value = "CRAMC"
$wsStock.Range("D9").Value = 317.2
$wsStock.Range("E9").Value = 2.3226
$wsStock.Range("H9").Value = 6325.5208
$wsStock.Range("B10").Value = "LGEINDIA"
$wsStock.Range("C10").Value = "LGEINDIA"
$wsStock.Range("D10").Value = 1633.4
$wsStock.Range("E10").Value = -0.946
$wsStock.Range("H10").Value = 110870.6825
$wsStock.Range("B11").Value = "TATACAP"
$wsStock.Range("C11").Value = "TATACAP"
$wsStock.Range("D11").Value = 329.3
$wsStock.Range("E11").Value = 0.1521
$wsStock.Range("H11").Value = 139783.5374
$wsStock.Range("B12").Value = "ELIQUID"
$wsStock.Range("C12").Value = "ELIQUID"
$wsStock.Range("D12").Value = 1004.85
$wsStock.Range("E12").Value = 0.0408
$wsStock.Range("H12").Value = 0
$wsStock.Range("B13").Value = "WEWORK"
$wsStock.Range("C13").Value = "WEWORK"
$wsStock.Range("D13").Value = 632.15
$wsStock.Range("E13").Value = -2.4008
$wsStock.Range("H13").Value = 8472.2803
$wsStock.Range("B14").Value = "GROWWRLTY"
$wsStock.Range("C14").Value = "GROWWRLTY"
$wsStock.Range("D14").Value = 10.8
$wsStock.Range("E14").Value = -0.4608
$wsStock.Range("H14").Value = 0
$wsStock.Range("B15").Value = "ADVANCE"
$wsStock.Range("C15").Value = "ADVANCE"
$wsStock.Range("D15").Value = 130.05
$wsStock.Range("E15").Value = -5.2666
$wsStock.Range("H15").Value = 836.0358
$wsStock.Range("B16").Value = "OMFREIGHT"
$wsStock.Range("C16").Value = "OMFREIGHT"
$wsStock.Range("D16").Value = 88.90000000000001
$wsStock.Range("E16").Value = -0.5926
$wsStock.Range("H16").Value = 299.3747
$wsStock.Range("B17").Value = "GLOTTIS"
$wsStock.Range("C17").Value = "GLOTTIS"
$wsStock.Range("D17").Value = 72.73999999999999
$wsStock.Range("E17").Value = -0.8587
$wsStock.Range("H17").Value = 672.1394
$wsStock.Range("B18").Value = "FABTECH"
$wsStock.Range("C18").Value = "FABTECH"
$wsStock.Range("D18").Value = 237.72
$wsStock.Range("E18").Value = 0.4734
$wsStock.Range("H18").Value = 1056.6843
$wsStock.Range("B19").Value = "PACEDIGITK"
$wsStock.Range("C19").Value = "PACEDIGITK"
$wsStock.Range("D19").Value = 218.85
$wsStock.Range("E19").Value = 0.1327
$wsStock.Range("H19").Value = 4723.9063
$wsStock.Range("B20").Value = "JAINREC"
$wsStock.Range("C20").Value = "JAINREC"
$wsStock.Range("D20").Value = 377.25
$wsStock.Range("E20").Value = 1.2208
$wsStock.Range("H20").Value = 13018.3623
$wsStock.Range("B21").Value = "EPACKPEB"
$wsStock.Range("C21").Value = "EPACKPEB"
$wsStock.Range("D21").Value = 301.45
$wsStock.Range("E21").Value = 1.979
$wsStock.Range("H21").Value = 3028.1254
$wsStock.Range("B22").Value = "BMWVENTLTD"
$wsStock.Range("C22").Value = "BMWVENTLTD"
$wsStock.Range("D22").Value = 69.25
$wsStock.Range("E22").Value = 0
$wsStock.Range("H22").Value = 600.5014
$wsStock.Range("B23").Value = "STYL"
$wsStock.Range("C23").Value = "STYL"
$wsStock.Range("D23").Value = 372.4
$wsStock.Range("E23").Value = -0.8388
$wsStock.Range("H23").Value = 6025.649
$wsStock.Range("B24").Value = "JARO"
$wsStock.Range("C24").Value = "JARO"
$wsStock.Range("D24").Value = 621.5
$wsStock.Range("E24").Value = -1.4821
$wsStock.Range("H24").Value = 1377.0134
$wsStock.Range("B25").Value = "SOLARWORLD"
$wsStock.Range("C25").Value = "SOLARWORLD"
$wsStock.Range("D25").Value = 309.1
$wsStock.Range("E25").Value = -0.6269
$wsStock.Range("H25").Value = 2679.0517
$wsStock.Range("B26").Value = "ARSSBL"
$wsStock.Range("C26").Value = "ARSSBL"
$wsStock.Range("D26").Value = 537.3
$wsStock.Range("E26").Value = 4.7266
$wsStock.Range("H26").Value = 3370.2277
$wsStock.Range("B27").Value = "GANESHCP"
$wsStock.Range("C27").Value = "GANESHCP"
$wsStock.Range("D27").Value = 274.4
$wsStock.Range("E27").Value = -2.7984
$wsStock.Range("H27").Value = 1108.9312
$wsStock.Range("B28").Value = "ATLANTAELE"
$wsStock.Range("C28").Value = "ATLANTAELE"
$wsStock.Range("D28").Value = 1003.05
$wsStock.Range("E28").Value = -1.7436
$wsStock.Range("H28").Value = 7713.116
$wsStock.Range("B29").Value = "GKENERGY"
$wsStock.Range("C29").Value = "GKENERGY"
$wsStock.Range("D29").Value = 213.85
$wsStock.Range("E29").Value = -0.7933
$wsStock.Range("H29").Value = 4337.2472
$wsStock.Range("B30").Value = "SAATVIKGL"
$wsStock.Range("C30").Value = "SAATVIKGL"
$wsStock.Range("D30").Value = 528.2
$wsStock.Range("E30").Value = -1.3079
$wsStock.Range("H30").Value = 6713.6863
$wsStock.Range("B31").Value = "IVALUE"
$wsStock.Range("C31").Value = "IVALUE"
$wsStock.Range("D31").Value = 281.45
$wsStock.Range("E31").Value = -0.3364
$wsStock.Range("H31").Value = 1506.8799
$wsStock.Range("B32").Value = "VMSTMT"
$wsStock.Range("C32").Value = "VMSTMT"
$wsStock.Range("D32").Value = 70.03
$wsStock.Range("E32").Value = -0.9056
$wsStock.Range("H32").Value = 347.5674
$wsStock.Range("B33").Value = "EUROPRATIK"
$wsStock.Range("C33").Value = "EUROPRATIK"
$wsStock.Range("D33").Value = 321.75
$wsStock.Range("E33").Value = 0.8147
$wsStock.Range("H33").Value = 3288.285
$wsStock.Range("B34").Value = "SHRINGARMS"
$wsStock.Range("C34").Value = "SHRINGARMS"
$wsStock.Range("D34").Value = 229.31
$wsStock.Range("E34").Value = -1.2616
$wsStock.Range("H34").Value = 2211.284
$wsStock.Range("B35").Value = "DEVX"
$wsStock.Range("C35").Value = "DEVX"
$wsStock.Range("D35").Value = 44.53
$wsStock.Range("E35").Value = -0.3803
$wsStock.Range("H35").Value = 401.605
$wsStock.Range("B36").Value = "URBANCO"
$wsStock.Range("C36").Value = "URBANCO"
$wsStock.Range("D36").Value = 148.9
$wsStock.Range("E36").Value = -2.0459
$wsStock.Range("H36").Value = 21380.5798
$wsStock.Range("B37").Value = "SML100CASE"
$wsStock.Range("C37").Value = "SML100CASE"
$wsStock.Range("D37").Value = 10.36
$wsStock.Range("E37").Value = -0.7663
$wsStock.Range("H37").Value = 0
$wsStock.Range("B38").Value = "AONEGOLD"
$wsStock.Range("C38").Value = "AONEGOLD"
$wsStock.Range("D38").Value = 11.28
$wsStock.Range("E38").Value = -0.2653
$wsStock.Range("H38").Value = 0
$wsStock.Range("B39").Value = "ELM250"
$wsStock.Range("C39").Value = "ELM250"
$wsStock.Range("D39").Value = 16.72
$wsStock.Range("E39").Value = 0.1797
$wsStock.Range("H39").Value = 0
$wsStock.Range("B40").Value = "AMANTA"
$wsStock.Range("C40").Value = "AMANTA"
$wsStock.Range("D40").Value = 122.52
$wsStock.Range("E40").Value = 1.407
$wsStock.Range("H40").Value = 475.7372
$wsStock.Range("B41").Value = "CPEDU"
$wsStock.Range("C41").Value = "CPEDU"
$wsStock.Range("D41").Value = 315.9
$wsStock.Range("E41").Value = 1.8539
$wsStock.Range("H41").Value = 574.7148999999999
$wsStock.Range("B42").Value = "AHCL"
$wsStock.Range("C42").Value = "AHCL"
$wsStock.Range("D42").Value = 139.27
$wsStock.Range("E42").Value = 3.1706
$wsStock.Range("H42").Value = 740.2409
$wsStock.Range("B43").Value = "STLNETWORK"
$wsStock.Range("C43").Value = "STLNETWORK"
$wsStock.Range("D43").Value = 26.59
$wsStock.Range("E43").Value = -0.412
$wsStock.Range("H43").Value = 1297.3822
$wsStock.Range("B44").Value = "VIKRAN"
$wsStock.Range("C44").Value = "VIKRAN"
$wsStock.Range("D44").Value = 98.05
$wsStock.Range("E44").Value = -1.783
$wsStock.Range("H44").Value = 2528.8166
$wsStock.Range("B45").Value = "MANUFGBEES"
$wsStock.Range("C45").Value = "MANUFGBEES"
$wsStock.Range("D45").Value = 151.77
$wsStock.Range("E45").Value = -1.011
$wsStock.Range("H45").Value = 0
$wsStock.Range("B46").Value = "MEIL"
$wsStock.Range("C46").Value = "MEIL"
$wsStock.Range("D46").Value = 461.15
$wsStock.Range("E46").Value = -0.7319
$wsStock.Range("H46").Value = 1274.1632
$wsStock.Range("B47").Value = "GROWWNXT50"
$wsStock.Range("C47").Value = "GROWWNXT50"
$wsStock.Range("D47").Value = 70.29000000000001
$wsStock.Range("E47").Value = -0.4109
$wsStock.Range("H47").Value = 0
$wsStock.Range("B48").Value = "SHREEJISPG"
$wsStock.Range("C48").Value = "SHREEJISPG"
$wsStock.Range("D48").Value = 270.05
$wsStock.Range("E48").Value = -0.7899
$wsStock.Range("H48").Value = 4399.6074
$wsStock.Range("B49").Value = "GEMAROMA"
$wsStock.Range("C49").Value = "GEMAROMA"
$wsStock.Range("D49").Value = 219.52
$wsStock.Range("E49").Value = -0.876
$wsStock.Range("H49").Value = 1146.7097
$wsStock.Range("B50").Value = "PATELRMART"
$wsStock.Range("C50").Value = "PATELRMART"
$wsStock.Range("D50").Value = 219.31
$wsStock.Range("E50").Value = -1.0646
$wsStock.Range("H50").Value = 732.5069999999999
$wsStock.Range("B51").Value = "VIKRAMSOLR"
$wsStock.Range("C51").Value = "VIKRAMSOLR"
$wsStock.Range("D51").Value = 322
$wsStock.Range("E51").Value = -1.5892
$wsStock.Range("H51").Value = 11647.2884
$wsStock.Range("B52").Value = "LTGILTCASE"
$wsStock.Range("C52").Value = "LTGILTCASE"
$wsStock.Range("D52").Value = 29.67
$wsStock.Range("E52").Value = 0.2365
$wsStock.Range("H52").Value = 0
$wsStock.Range("B53").Value = "REGAAL"
$wsStock.Range("C53").Value = "REGAAL"
$wsStock.Range("D53").Value = 89.13
$wsStock.Range("E53").Value = -0.8675
$wsStock.Range("H53").Value = 915.5742
$wsStock.Range("B54").Value = "BLUESTONE"
$wsStock.Range("C54").Value = "BLUESTONE"
$wsStock.Range("D54").Value = 711.95
$wsStock.Range("E54").Value = 0.1266
$wsStock.Range("H54").Value = 10773.2539
$wsStock.Range("B55").Value = "MOSILVER"
$wsStock.Range("C55").Value = "MOSILVER"
$wsStock.Range("D55").Value = 145.9
$wsStock.Range("E55").Value = -1.5054
$wsStock.Range("H55").Value = 0
$wsStock.Range("B56").Value = "ALLTIME"
$wsStock.Range("C56").Value = "ALLTIME"
$wsStock.Range("D56").Value = 308.75
$wsStock.Range("E56").Value = 2.66
$wsStock.Range("H56").Value = 2022.5526
$wsStock.Range("B57").Value = "JSWCEMENT"
$wsStock.Range("C57").Value = "JSWCEMENT"
$wsStock.Range("D57").Value = 134.98
$wsStock.Range("E57").Value = -0.4793
$wsStock.Range("H57").Value = 18402.6999
$wsStock.Range("B58").Value = "SBILIQETF"
$wsStock.Range("C58").Value = "SBILIQETF"
$wsStock.Range("D58").Value = 1012.94
$wsStock.Range("E58").Value = 0.0296
$wsStock.Range("H58").Value = 0
$wsStock.Range("B59").Value = "HILINFRA"
$wsStock.Range("C59").Value = "HILINFRA"
$wsStock.Range("D59").Value = 77.23
$wsStock.Range("E59").Value = -0.3998
$wsStock.Range("H59").Value = 0
$wsStock.Range("B60").Value = "GROWWPOWER"
$wsStock.Range("C60").Value = "GROWWPOWER"
$wsStock.Range("D60").Value = 10.28
$wsStock.Range("E60").Value = -0.9634
$wsStock.Range("H60").Value = 0
$wsStock.Range("B61").Value = "LOTUSDEV"
$wsStock.Range("C61").Value = "LOTUSDEV"
$wsStock.Range("D61").Value = 177.82
$wsStock.Range("E61").Value = 0.3669
$wsStock.Range("H61").Value = 8690.485000000001
$wsStock.Range("B62").Value = "MBEL"
$wsStock.Range("C62").Value = "MBEL"
$wsStock.Range("D62").Value = 450.2
$wsStock.Range("E62").Value = -0.7714
$wsStock.Range("H62").Value = 2572.8126
$wsStock.Range("B63").Value = "LAXMIINDIA"
$wsStock.Range("C63").Value = "LAXMIINDIA"
$wsStock.Range("D63").Value = 145.62
$wsStock.Range("E63").Value = -1.1942
$wsStock.Range("H63").Value = 761.1248000000001
$wsStock.Range("B64").Value = "CPPLUS"
$wsStock.Range("C64").Value = "CPPLUS"
$wsStock.Range("D64").Value = 1322.1
$wsStock.Range("E64").Value = -0.264
$wsStock.Range("H64").Value = 15497.9053
$wsStock.Range("B65").Value = "SHANTIGOLD"
$wsStock.Range("C65").Value = "SHANTIGOLD"
$wsStock.Range("D65").Value = 241.57
$wsStock.Range("E65").Value = -1.6409
$wsStock.Range("H65").Value = 1741.6231
$wsStock.Range("B66").Value = "MOGOLD"
$wsStock.Range("C66").Value = "MOGOLD"
$wsStock.Range("D66").Value = 119.65
$wsStock.Range("E66").Value = -0.5403
$wsStock.Range("H66").Value = 0
$wsStock.Range("B67").Value = "BRIGHOTEL"
$wsStock.Range("C67").Value = "BRIGHOTEL"
$wsStock.Range("D67").Value = 82.39
$wsStock.Range("E67").Value = -0.9855
$wsStock.Range("H67").Value = 3129.5229
$wsStock.Range("B68").Value = "INDIQUBE"
$wsStock.Range("C68").Value = "INDIQUBE"
$wsStock.Range("D68").Value = 212.64
$wsStock.Range("E68").Value = -0.7561
$wsStock.Range("H68").Value = 4465.6847
$wsStock.Range("B69").Value = "EBGNG"
$wsStock.Range("C69").Value = "EBGNG"
$wsStock.Range("D69").Value = 346.65
$wsStock.Range("E69").Value = 3.2311
$wsStock.Range("H69").Value = 3952.2092
$wsStock.Range("B70").Value = "LIQGRWBEES"
$wsStock.Range("C70").Value = "LIQGRWBEES"
$wsStock.Range("D70").Value = 1014.74
$wsStock.Range("E70").Value = 0.0246
$wsStock.Range("H70").Value = 0
$wsStock.Range("B71").Value = "CHEMBONDCH"
$wsStock.Range("C71").Value = "CHEMBONDCH"
$wsStock.Range("D71").Value = 153.35
$wsStock.Range("E71").Value = -1.6987
$wsStock.Range("H71").Value = 412.459
$wsStock.Range("B72").Value = "GROWWNIFTY"
$wsStock.Range("C72").Value = "GROWWNIFTY"
$wsStock.Range("D72").Value = 10.29
$wsStock.Range("E72").Value = -0.3872
$wsStock.Range("H72").Value = 0
$wsStock.Range("B73").Value = "ANTHEM"
$wsStock.Range("C73").Value = "ANTHEM"
$wsStock.Range("D73").Value = 702.25
$wsStock.Range("E73").Value = -0.1209
$wsStock.Range("H73").Value = 39439.0658
$wsStock.Range("B74").Value = "QUALITY30"
$wsStock.Range("C74").Value = "QUALITY30"
$wsStock.Range("D74").Value = 21.05
$wsStock.Range("E74").Value = -0.8945
$wsStock.Range("H74").Value = 0
$wsStock.Range("B75").Value = "SMARTWORKS"
$wsStock.Range("C75").Value = "SMARTWORKS"
$wsStock.Range("D75").Value = 606.65
$wsStock.Range("E75").Value = 2.0867
$wsStock.Range("H75").Value = 6931.2448
$wsStock.Range("B76").Value = "TRAVELFOOD"
$wsStock.Range("C76").Value = "TRAVELFOOD"
$wsStock.Range("D76").Value = 1316.3
$wsStock.Range("E76").Value = 0.1141
$wsStock.Range("H76").Value = 17332.9705
